$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "ECs" target-cluster row (originally row 2); remaining rows shift up.
$ws.Rows.Item(2).Delete()

# Row 2 (was row 3, Target=FAPs) gets refreshed receptor/edge specificity values
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.8326193333333333
$ws.Range("N2").Value = 2.497858
$ws.Range("O2").Value = 0.9388124812781204
$ws.Range("P2").Value = 0.9388124812781203
$ws.Range("Q2").Value = 0.1814754895751111
$ws.Range("R2").Value = 1.633279406176
$ws.Range("S2").Value = 0.9388124812781204
$ws.Range("T2").Value = 0.9388124812781203

# Row 3 (was row 4, Target=MuSCs) gets refreshed receptor/edge values from new TPM data
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.05426633333333333
$ws.Range("N3").Value = 0.162799
$ws.Range("O3").Value = 0.0611875187218796
$ws.Range("P3").Value = 0.06118751872187959
$ws.Range("Q3").Value = 0.01182774530311111
$ws.Range("R3").Value = 0.106449707728
$ws.Range("S3").Value = 0.0611875187218796
$ws.Range("T3").Value = 0.06118751872187959
